$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 3-4: refreshed siniestro numbers (style / formatting unchanged) --
# These look like numbers but must stay as text; the leading apostrophe
# forces text interpretation (the cells already carry the workbook's
# quotePrefix text style, so we only touch the value).
$ws.Range("F3").Value = "'1120194100448 "
$ws.Range("F4").Value = "'1220194200694  "

# --- Rows 5-7: "Anular" action test cases -----------------------------------
# These cells previously had the shaded/bordered quotePrefix style; the
# corrected test data drops that shading, so clear formatting first and then
# write the (text) value.
$ws.Range("F5").ClearFormats()
$ws.Range("F5").Value = "'0420172010228 "
$ws.Range("G5").Value = "Anular"

$ws.Range("F6").ClearFormats()
$ws.Range("F6").Value = "'1120170200969 "
$ws.Range("G6").Value = "Anular"

$ws.Range("F7").ClearFormats()
$ws.Range("F7").Value = "'1220170301466 "
$ws.Range("G7").Value = "Anular"

# --- Rows 8-10: "Obtener numero" action test cases --------------------------
$ws.Range("F8").ClearFormats()
$ws.Range("F8").Value = "'0420172010228 "
$ws.Range("G8").Value = "Obtener número"

$ws.Range("B9").Value = $ws.Range("B8").Value2
$ws.Range("C9").Value = $ws.Range("C8").Value2
$ws.Range("D9").Value = $ws.Range("D8").Value2
$ws.Range("E9").Value = $ws.Range("E8").Value2
$ws.Range("F9").Value = "'1120170200969 "
$ws.Range("G9").Value = "Obtener número"

$ws.Range("B10").Value = $ws.Range("B8").Value2
$ws.Range("C10").Value = $ws.Range("C8").Value2
$ws.Range("D10").Value = $ws.Range("D8").Value2
$ws.Range("E10").Value = $ws.Range("E8").Value2
$ws.Range("F10").Value = "'1220170301466 "
$ws.Range("G10").Value = "Obtener número"

# --- Restore the active selection to F4, as recorded in the saved workbook -
$ws.Range("F4").Select() | Out-Null
